# Auto-generated edit script: updates market/profit columns (H-N) per the commit diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 774.5789
$ws.Range("I19").Value = 626.61536
$ws.Range("J19").Value = 1095.1666
$ws.Range("K19").Value = 626.61536
$ws.Range("L19").Value = 1095.1666
$ws.Range("M19").Value = -451.61536
$ws.Range("N19").Value = -1445.1666
$ws.Range("H32").Value = 1473.5555
$ws.Range("J32").Value = 1664.5
$ws.Range("L32").Value = 1664.5
$ws.Range("N32").Value = -2316.5
$ws.Range("H41").Value = 1053.3846
$ws.Range("I41").Value = 2509.4
$ws.Range("K41").Value = 2509.4
$ws.Range("M41").Value = -2069.4
$ws.Range("H64").Value = 5099.294
$ws.Range("J64").Value = 7444.6665
$ws.Range("L64").Value = 7444.6665
$ws.Range("N64").Value = -7940.6665
$ws.Range("H67").Value = 5099.294
$ws.Range("J67").Value = 7444.6665
$ws.Range("L67").Value = 7444.6665
$ws.Range("N67").Value = -9160.666499999999
$ws.Range("H76").Value = 142860700
$ws.Range("J76").Value = 4159.4
$ws.Range("L76").Value = 4159.4
$ws.Range("N76").Value = -4789.4
$ws.Range("H79").Value = 142860700
$ws.Range("J79").Value = 4159.4
$ws.Range("L79").Value = 4159.4
$ws.Range("N79").Value = -6343.4
$ws.Range("H98").Value = 4864.75
$ws.Range("I98").Value = 1431.3077
$ws.Range("J98").Value = 49499.5
$ws.Range("K98").Value = 1431.3077
$ws.Range("L98").Value = 49499.5
$ws.Range("M98").Value = 66.69229999999993
$ws.Range("N98").Value = -52495.5
$ws.Range("H113").Value = 6710.45
$ws.Range("I113").Value = 6519
$ws.Range("K113").Value = 6519
$ws.Range("M113").Value = -3265
$ws.Range("H122").Value = 4864.75
$ws.Range("I122").Value = 1431.3077
$ws.Range("J122").Value = 49499.5
$ws.Range("K122").Value = 4293.9231
$ws.Range("L122").Value = 148498.5
$ws.Range("M122").Value = -1843.9231
$ws.Range("N122").Value = -153398.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 18331.666
$ws.Range("I11").Value = 15000
$ws.Range("K11").Value = 15000
$ws.Range("M11").Value = -14856
$ws.Range("H32").Value = 4012.25
$ws.Range("I32").Value = 3156.348
$ws.Range("J32").Value = 13855.125
$ws.Range("K32").Value = 3156.348
$ws.Range("L32").Value = 13855.125
$ws.Range("M32").Value = -2869.348
$ws.Range("N32").Value = -14429.125
$ws.Range("H63").Value = 1548.6666
$ws.Range("I63").Value = 1673
$ws.Range("J63").Value = 1300
$ws.Range("K63").Value = 1673
$ws.Range("L63").Value = 1300
$ws.Range("M63").Value = -987
$ws.Range("N63").Value = -2672
$ws.Range("H66").Value = 1548.6666
$ws.Range("I66").Value = 1673
$ws.Range("J66").Value = 1300
$ws.Range("K66").Value = 8365
$ws.Range("L66").Value = 6500
$ws.Range("M66").Value = -4933
$ws.Range("N66").Value = -13364
$ws.Range("H74").Value = 5191.778
$ws.Range("I74").Value = 3499
$ws.Range("J74").Value = 5675.4287
$ws.Range("K74").Value = 3499
$ws.Range("L74").Value = 5675.4287
$ws.Range("M74").Value = -2625
$ws.Range("N74").Value = -7423.4287
$ws.Range("H77").Value = 5191.778
$ws.Range("I77").Value = 3499
$ws.Range("J77").Value = 5675.4287
$ws.Range("K77").Value = 17495
$ws.Range("L77").Value = 28377.1435
$ws.Range("M77").Value = -13127
$ws.Range("N77").Value = -37113.14350000001
$ws.Range("H88").Value = 9811.666999999999
$ws.Range("J88").Value = 12605.111
$ws.Range("L88").Value = 12605.111
$ws.Range("N88").Value = -13417.111
$ws.Range("H91").Value = 9811.666999999999
$ws.Range("J91").Value = 12605.111
$ws.Range("L91").Value = 12605.111
$ws.Range("N91").Value = -15413.111
$ws.Range("H97").Value = 3208.389
$ws.Range("I97").Value = 664.95654
$ws.Range("J97").Value = 7708.3076
$ws.Range("K97").Value = 664.95654
$ws.Range("L97").Value = 7708.3076
$ws.Range("M97").Value = -168.95654
$ws.Range("N97").Value = -8700.3076

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2200.6924
$ws.Range("J86").Value = 2606.4
$ws.Range("L86").Value = 2606.4
$ws.Range("N86").Value = -4852.4
$ws.Range("H89").Value = 2200.6924
$ws.Range("J89").Value = 2606.4
$ws.Range("L89").Value = 13032
$ws.Range("N89").Value = -24264
$ws.Range("H105").Value = 3751.75
$ws.Range("I105").Value = 3204.3333
$ws.Range("J105").Value = 4199.636
$ws.Range("K105").Value = 3204.3333
$ws.Range("L105").Value = 4199.636
$ws.Range("M105").Value = -1457.3333
$ws.Range("N105").Value = -7693.636

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 379.2857
$ws.Range("I16").Value = 379.2857
$ws.Range("K16").Value = 379.2857
$ws.Range("M16").Value = -92.28570000000002
$ws.Range("H38").Value = 10000
$ws.Range("I38").Value = 10000
$ws.Range("K38").Value = 10000
$ws.Range("M38").Value = -9623
$ws.Range("H46").Value = 10000
$ws.Range("I46").Value = 10000
$ws.Range("K46").Value = 10000
$ws.Range("M46").Value = -9789
$ws.Range("H62").Value = 3614.077
$ws.Range("J62").Value = 4198
$ws.Range("L62").Value = 4198
$ws.Range("N62").Value = -5446
$ws.Range("H65").Value = 3614.077
$ws.Range("J65").Value = 4198
$ws.Range("L65").Value = 20990
$ws.Range("N65").Value = -27230
$ws.Range("H99").Value = 2697.5386
$ws.Range("I99").Value = 2439.4285
$ws.Range("J99").Value = 2998.6667
$ws.Range("K99").Value = 2439.4285
$ws.Range("L99").Value = 2998.6667
$ws.Range("M99").Value = -941.4285
$ws.Range("N99").Value = -5994.6667
$ws.Range("H113").Value = 379.2857
$ws.Range("I113").Value = 379.2857
$ws.Range("K113").Value = 379.2857
$ws.Range("M113").Value = 1790.7143
$ws.Range("H126").Value = 2697.5386
$ws.Range("I126").Value = 2439.4285
$ws.Range("J126").Value = 2998.6667
$ws.Range("K126").Value = 7318.2855
$ws.Range("L126").Value = 8996.000100000001
$ws.Range("M126").Value = -4848.2855
$ws.Range("N126").Value = -13936.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 341.16666
$ws.Range("J23").Value = 386.14285
$ws.Range("L23").Value = 1158.42855
$ws.Range("N23").Value = -1628.42855
$ws.Range("H40").Value = 53
$ws.Range("I40").Value = 53.25
$ws.Range("K40").Value = 213
$ws.Range("M40").Value = -144
$ws.Range("H131").Value = 9663226
$ws.Range("J131").Value = 9260744
$ws.Range("L131").Value = 27782232
$ws.Range("N131").Value = -27792312
$ws.Range("H137").Value = 44157.816
$ws.Range("I137").Value = 96861.17999999999
$ws.Range("K137").Value = 290583.54
$ws.Range("M137").Value = -285483.54

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 11374.5
$ws.Range("H44").Value = 24974.5
$ws.Range("I44").Value = 9950
$ws.Range("K44").Value = 9950
$ws.Range("M44").Value = -9354
$ws.Range("H80").Value = 4149.2
$ws.Range("I80").Value = 4082.3333
$ws.Range("J80").Value = 4249.5
$ws.Range("K80").Value = 4082.3333
$ws.Range("L80").Value = 4249.5
$ws.Range("M80").Value = -3084.3333
$ws.Range("N80").Value = -6245.5
$ws.Range("H83").Value = 4149.2
$ws.Range("I83").Value = 4082.3333
$ws.Range("J83").Value = 4249.5
$ws.Range("K83").Value = 20411.6665
$ws.Range("L83").Value = 21247.5
$ws.Range("M83").Value = -15419.6665
$ws.Range("N83").Value = -31231.5
$ws.Range("H102").Value = 23767.25
$ws.Range("I102").Value = 25609.727
$ws.Range("K102").Value = 25609.727
$ws.Range("M102").Value = -23987.727
$ws.Range("H122").Value = 65554.82000000001
$ws.Range("I122").Value = 131019.375
$ws.Range("J122").Value = 7364.1113
$ws.Range("K122").Value = 393058.125
$ws.Range("L122").Value = 22092.3339
$ws.Range("M122").Value = -390608.125
$ws.Range("N122").Value = -26992.3339

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10000
$ws.Range("I7").Value = 10000
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 10000
$ws.Range("L7").Value = 0
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -9888
$ws.Range("H40").Value = 3886.5
$ws.Range("I40").Value = 2649
$ws.Range("K40").Value = 2649
$ws.Range("M40").Value = -2513
$ws.Range("H55").Value = 961.375
$ws.Range("I55").Value = 218.63637
$ws.Range("J55").Value = 1589.8462
$ws.Range("K55").Value = 218.63637
$ws.Range("L55").Value = 1589.8462
$ws.Range("M55").Value = -45.63637
$ws.Range("N55").Value = -1935.8462
$ws.Range("H126").Value = 10000
$ws.Range("I126").Value = 10000
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 30000
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -27530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 23166.334
$ws.Range("I21").Value = 21332.666
$ws.Range("J21").Value = 25000
$ws.Range("K21").Value = 21332.666
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = -21097.666
$ws.Range("N21").Value = -25470
$ws.Range("H35").Value = 23166.334
$ws.Range("I35").Value = 21332.666
$ws.Range("J35").Value = 25000
$ws.Range("K35").Value = 21332.666
$ws.Range("L35").Value = 25000
$ws.Range("M35").Value = -21042.666
$ws.Range("N35").Value = -25580
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -51137
